$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.973.92'
$ws.Range("E2").Value = '  +0.86%  '
$ws.Range("D3").Value = '2.310.76'
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = "'540.07"
$ws.Range("E5").Value = '  -1.31%  '
$ws.Range("D6").Value = "'131.21"
$ws.Range("E6").Value = '  -0.66%  '
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("D8").Value = "'0.582"
$ws.Range("E8").Value = '  +1.97%  '
$ws.Range("D9").Value = '2.305.46'
$ws.Range("E9").Value = '  +0.17%  '
$ws.Range("D10").Value = "'0.0998"
$ws.Range("E10").Value = '  -1.60%  '
$ws.Range("D11").Value = "'5.48"
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("E12").Value = '  +0.43%  '
$ws.Range("E13").Value = '  -0.13%  '
$ws.Range("D14").Value = "'23.67"
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("D15").Value = '2.722.41'
$ws.Range("E15").Value = '  +0.50%  '
$ws.Range("D16").Value = '58.925.57'
$ws.Range("E16").Value = '  +0.95%  '
$ws.Range("E17").Value = '  -0.24%  '
$ws.Range("D18").Value = '2.312.29'
$ws.Range("E18").Value = '  +1.48%  '
$ws.Range("D19").Value = "'10.53"
$ws.Range("E19").Value = '  -0.78%  '
$ws.Range("E20").Value = '  -2.69%  '
$ws.Range("D21").Value = "'312.89"
$ws.Range("E21").Value = '  -0.21%  '
$ws.Range("D22").Value = "'6.59"
$ws.Range("E22").Value = '  +2.09%  '
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").Value = "'62.58"
$ws.Range("E24").Value = '  -0.86%  '
$ws.Range("E25").Value = '  +2.70%  '
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("D27").Value = "'7.91"
$ws.Range("E27").Value = '  -1.41%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("D29").Value = "'171.23"
$ws.Range("E29").Value = '  +0.48%  '
$ws.Range("D30").Value = "'1.71"
$ws.Range("E30").Value = '  -1.99%  '
$ws.Range("E31").Value = '  +6.83%  '
$ws.Range("D32").Value = '0.0₃0733'
$ws.Range("E32").Value = '  +1.93%  '
$ws.Range("D33").Value = "'5.85"
$ws.Range("E33").Value = '  +1.76%  '
$ws.Range("E34").Value = '  +15.07%  '
$ws.Range("D35").Value = "'0.383"
$ws.Range("E35").Value = '  +0.79%  '
$ws.Range("D37").Value = "'17.79"
$ws.Range("E37").Value = '  +0.24%  '
$ws.Range("E38").Value = '  +0.13%  '
$ws.Range("D39").Value = "'4.03"
$ws.Range("E39").Value = '  +2.52%  '
$ws.Range("D40").Value = "'312.35"
$ws.Range("E40").Value = '  +7.15%  '
$ws.Range("D41").Value = "'37.88"
$ws.Range("E41").Value = '  -0.55%  '
$ws.Range("E42").Value = '  +0.61%  '
$ws.Range("D43").Value = "'141.63"
$ws.Range("E43").Value = '  +0.70%  '
$ws.Range("E44").Value = '  +0.14%  '
$ws.Range("D45").Value = "'0.0951"
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("E46").Value = '  -1.84%  '
$ws.Range("D47").Value = "'0.556"
$ws.Range("E47").Value = '  +0.28%  '
$ws.Range("D48").Value = "'18.19"
$ws.Range("E48").Value = '  -0.65%  '
$ws.Range("E49").Value = '  -2.38%  '
$ws.Range("D50").Value = "'11.01"
$ws.Range("E51").Value = '  -0.19%  '
